$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 2422.2
$ws.Range("I46").Value = 2501.2727
$ws.Range("J46").Value = 2204.75
$ws.Range("K46").Value = 7503.8181
$ws.Range("L46").Value = 6614.25
$ws.Range("M46").Value = -7384.8181
$ws.Range("N46").Value = -6852.25

$ws.Range("H60").Value = 2422.2
$ws.Range("I60").Value = 2501.2727
$ws.Range("J60").Value = 2204.75
$ws.Range("K60").Value = 7503.8181
$ws.Range("L60").Value = 6614.25
$ws.Range("M60").Value = -7019.8181
$ws.Range("N60").Value = -7582.25

$ws.Range("H116").Value = 3967.9614
$ws.Range("I116").Value = 3489.3635
$ws.Range("J116").Value = 4318.933
$ws.Range("K116").Value = 3489.3635
$ws.Range("L116").Value = 4318.933
$ws.Range("M116").Value = -47.36349999999993
$ws.Range("N116").Value = -11202.933

$ws.Range("H132").Value = 16166.719
$ws.Range("I132").Value = 2207.5576
$ws.Range("K132").Value = 6622.6728
$ws.Range("M132").Value = -4092.6728

$ws.Range("H137").Value = 12015.429
$ws.Range("I137").Value = 9730.143
$ws.Range("J137").Value = 16586
$ws.Range("K137").Value = 29190.429
$ws.Range("L137").Value = 49758
$ws.Range("M137").Value = -26640.429
$ws.Range("N137").Value = -54858

$ws.Range("H138").Value = 2693.0642
$ws.Range("I138").Value = 1482.7858
$ws.Range("J138").Value = 2957.8125
$ws.Range("K138").Value = 4448.357400000001
$ws.Range("L138").Value = 8873.4375
$ws.Range("M138").Value = 691.6425999999992
$ws.Range("N138").Value = -19153.4375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1910.1666
$ws.Range("I2").Value = 1542.1666
$ws.Range("J2").Value = 3382.1667
$ws.Range("K2").Value = 1542.1666
$ws.Range("L2").Value = 3382.1667
$ws.Range("M2").Value = -1429.1666
$ws.Range("N2").Value = -3608.1667

$ws.Range("H61").Value = 2750.8823
$ws.Range("I61").Value = 2034.5555
$ws.Range("K61").Value = 2034.5555
$ws.Range("M61").Value = -1822.5555

$ws.Range("H116").Value = 1910.1666
$ws.Range("I116").Value = 1542.1666
$ws.Range("J116").Value = 3382.1667
$ws.Range("K116").Value = 1542.1666
$ws.Range("L116").Value = 3382.1667
$ws.Range("M116").Value = 751.8334
$ws.Range("N116").Value = -7970.1667

$ws.Range("H132").Value = 13891936
$ws.Range("I132").Value = 29413932
$ws.Range("J132").Value = 3833.842
$ws.Range("K132").Value = 88241796
$ws.Range("L132").Value = 11501.526
$ws.Range("M132").Value = -88239266
$ws.Range("N132").Value = -16561.526

$ws.Range("H136").Value = 2750.8823
$ws.Range("I136").Value = 2034.5555
$ws.Range("K136").Value = 6103.666499999999
$ws.Range("M136").Value = -3553.666499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1910.1666
$ws.Range("I3").Value = 1542.1666
$ws.Range("J3").Value = 3382.1667
$ws.Range("K3").Value = 1542.1666
$ws.Range("L3").Value = 3382.1667
$ws.Range("M3").Value = -1428.1666
$ws.Range("N3").Value = -3610.1667

$ws.Range("H96").Value = 14126.087
$ws.Range("I96").Value = 3975
$ws.Range("J96").Value = 16263.158
$ws.Range("K96").Value = 3975
$ws.Range("L96").Value = 16263.158
$ws.Range("M96").Value = -1229
$ws.Range("N96").Value = -21755.158

$ws.Range("H133").Value = 53995
$ws.Range("J133").Value = 53995
$ws.Range("L133").Value = 53995
$ws.Range("N133").Value = -64115

$ws.Range("H134").Value = 2913.0386
$ws.Range("I134").Value = 2434.3684
$ws.Range("K134").Value = 7303.1052
$ws.Range("M134").Value = -4768.1052

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 8777356
$ws.Range("I31").Value = 2960.5833
$ws.Range("J31").Value = 12827076
$ws.Range("K31").Value = 2960.5833
$ws.Range("L31").Value = 12827076
$ws.Range("M31").Value = -2665.5833
$ws.Range("N31").Value = -12827666

$ws.Range("H34").Value = 8777356
$ws.Range("I34").Value = 2960.5833
$ws.Range("J34").Value = 12827076
$ws.Range("K34").Value = 2960.5833
$ws.Range("L34").Value = 12827076
$ws.Range("M34").Value = -2758.5833
$ws.Range("N34").Value = -12827480

$ws.Range("H58").Value = 15153244
$ws.Range("I58").Value = 973.44446
$ws.Range("J58").Value = 33335968
$ws.Range("K58").Value = 973.44446
$ws.Range("L58").Value = 33335968
$ws.Range("M58").Value = -770.44446
$ws.Range("N58").Value = -33336374

$ws.Range("H99").Value = 2158.111
$ws.Range("I99").Value = 2002.875
$ws.Range("J99").Value = 3400
$ws.Range("K99").Value = 2002.875
$ws.Range("L99").Value = 3400
$ws.Range("M99").Value = -504.875
$ws.Range("N99").Value = -6396

$ws.Range("H126").Value = 2158.111
$ws.Range("I126").Value = 2002.875
$ws.Range("J126").Value = 3400
$ws.Range("K126").Value = 6008.625
$ws.Range("L126").Value = 10200
$ws.Range("M126").Value = -3538.625
$ws.Range("N126").Value = -15140

$ws.Range("H132").Value = 25034.033
$ws.Range("I132").Value = 1342.9574
$ws.Range("J132").Value = 104568.36
$ws.Range("K132").Value = 4028.8722
$ws.Range("L132").Value = 313705.08
$ws.Range("M132").Value = -1498.8722
$ws.Range("N132").Value = -318765.08

$ws.Range("H136").Value = 15153244
$ws.Range("I136").Value = 973.44446
$ws.Range("J136").Value = 33335968
$ws.Range("K136").Value = 2920.33338
$ws.Range("L136").Value = 100007904
$ws.Range("M136").Value = -370.33338
$ws.Range("N136").Value = -100013004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 52.357143
$ws.Range("I14").Value = 52.357143
$ws.Range("K14").Value = 157.071429
$ws.Range("M14").Value = 15.92857100000001

$ws.Range("H107").Value = 3244.7837
$ws.Range("I107").Value = 5836.3335
$ws.Range("J107").Value = 2411.7856
$ws.Range("K107").Value = 17509.0005
$ws.Range("L107").Value = 7235.3568
$ws.Range("M107").Value = -15589.0005
$ws.Range("N107").Value = -11075.3568

$ws.Range("H131").Value = 783.3200000000001
$ws.Range("I131").Value = 381.36365
$ws.Range("J131").Value = 833
$ws.Range("K131").Value = 1144.09095
$ws.Range("L131").Value = 2499
$ws.Range("M131").Value = 3895.90905
$ws.Range("N131").Value = -12579

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 55663.668
$ws.Range("J133").Value = 55663.668
$ws.Range("L133").Value = 55663.668
$ws.Range("N133").Value = -65783.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3067.5293
$ws.Range("I7").Value = 2528
$ws.Range("J7").Value = 4056.6667
$ws.Range("K7").Value = 2528
$ws.Range("L7").Value = 4056.6667
$ws.Range("M7").Value = -2416
$ws.Range("N7").Value = -4280.6667

$ws.Range("H40").Value = 6166.1665
$ws.Range("I40").Value = 4998.5
$ws.Range("J40").Value = 6750
$ws.Range("K40").Value = 4998.5
$ws.Range("L40").Value = 6750
$ws.Range("M40").Value = -4862.5
$ws.Range("N40").Value = -7022

$ws.Range("H126").Value = 3067.5293
$ws.Range("I126").Value = 2528
$ws.Range("J126").Value = 4056.6667
$ws.Range("K126").Value = 7584
$ws.Range("L126").Value = 12170.0001
$ws.Range("M126").Value = -5114
$ws.Range("N126").Value = -17110.0001

$ws.Range("H132").Value = 2967.0195
$ws.Range("I132").Value = 2131.4243
$ws.Range("K132").Value = 6394.2729
$ws.Range("M132").Value = -3864.2729

$ws.Range("H136").Value = 1740.8379
$ws.Range("I136").Value = 1188.6666
$ws.Range("J136").Value = 2760.2307
$ws.Range("K136").Value = 3565.9998
$ws.Range("L136").Value = 8280.6921
$ws.Range("M136").Value = -1015.9998
$ws.Range("N136").Value = -13380.6921

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2484
$ws.Range("I126").Value = 2380.8
$ws.Range("K126").Value = 7142.400000000001
$ws.Range("M126").Value = -4672.400000000001

$ws.Range("H132").Value = 14287648
$ws.Range("I132").Value = 1547.6957
$ws.Range("J132").Value = 41669340
$ws.Range("K132").Value = 4643.0871
$ws.Range("L132").Value = 125008020
$ws.Range("M132").Value = -2113.0871
$ws.Range("N132").Value = -125013080

$ws.Range("H136").Value = 304102.44
$ws.Range("I136").Value = 357985.06
$ws.Range("J136").Value = 2359.6
$ws.Range("K136").Value = 1073955.18
$ws.Range("L136").Value = 7078.799999999999
$ws.Range("M136").Value = -1071405.18
$ws.Range("N136").Value = -12178.8
